$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recipe")
$ws.Activate()

# Clear the cocktail_pk values in column B for rows 4-17 (formatting stays column-driven)
$ws.Range("B4:B17").Clear()

# Clear the cocktail_kor text in column C for row 2 and rows 11-25 (keep cell formatting)
$ws.Range("C2").ClearContents()
$ws.Range("C11:C25").ClearContents()

# Update the active selection shown in the sheet view
$ws.Range("E10").Select()
